# Automatische test-sync: 2025-08-05 17:18:50
# Appends a new log row (row 16) to the "Logs" sheet and updates the
# "Dashboard" summary count for "Planning / Afspraak" from 9 to 10.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of data to append at row 16 on the Logs sheet.
$newRow = @(
    "Wil je dit oppakken?",
    "mailmind.test@zohomail.eu",
    "Testmail #2: Wil je dit oppakken?",
    "Planning / Afspraak",
    "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl.",
    "2025-08-05 17:17:52",
    "Ja",
    "Ja",
    "Nee",
    "Nee"
)

$rowIndex = 16
for ($col = 1; $col -le $newRow.Length; $col++) {
    $logs.Cells.Item($rowIndex, $col).Value = $newRow[$col - 1]
}

# Update the Dashboard count for "Planning / Afspraak" (B2: 9 -> 10).
$dashboard.Range("B2").Value = 10

# Expand the conditional formatting ranges so they cover the new row 16
# (D2:D15 -> D2:D16, G2:G15 -> G2:G16, H2:H15 -> H2:H16, I2:I15 -> I2:I16,
#  J2:J15 -> J2:J16).
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$col`2:$col`15")
    $newRange = $logs.Range("$col`2:$col`16")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
